# 🔄 MAJ automatique BRVM via GitHub Actions
#
# This script re-applies the daily BRVM "Recommandations" / "Top_YTD" refresh:
#  - updates the aggregate sector rows (2-15) with the new "Variation Totale (%)"
#    and "Derniere Variation (%)" figures,
#  - re-ranks the individual-title rows (22-38) by the new "Variation Totale (%)",
#    rewriting each row's title/day-counts/variations/recommendation/strategy,
#  - drops the 4 titles that fell out of the ranked list (rows 39-42),
#  - refreshes the "Top_YTD" sheet's YTD progression figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Recommandations"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Recommandations")

# Rows 2-15: sector aggregate rows - only columns D (Variation Totale) and
# E (Derniere Variation) change; title/B/C/F/G stay as-is.
$sectorUpdates = @(
    @{ Row = 2;  D = 2535.03;  E = 109.62 },
    @{ Row = 3;  D = 2025;     E = 690 },
    @{ Row = 4;  D = 1985;     E = 665 },
    @{ Row = 5;  D = 1826.34;  E = 611.01 },
    @{ Row = 6;  D = 1451.46;  E = 494.9 },
    @{ Row = 7;  D = 1138.77;  E = 375.36 },
    @{ Row = 8;  D = 1089.23;  E = 362.26 },
    @{ Row = 9;  D = 491.28;   E = 168.85 },
    @{ Row = 10; D = 413.32;   E = 138.95 },
    @{ Row = 11; D = 408.94;   E = 138.18 },
    @{ Row = 12; D = 401.9;    E = 135.8 },
    @{ Row = 13; D = 388.19;   E = 126.51 },
    @{ Row = 14; D = 332.46;   E = 111.9 },
    @{ Row = 15; D = 288.61;   E = 96 }
)

foreach ($u in $sectorUpdates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Rows 16-21 are unchanged.

# Rows 22-38: individual-title ranking, fully rewritten (title, hausse/baisse
# day counts, variations, recommandation, strategie).
$titleUpdates = @(
    @{ Row = 22; A = "SICABLE CI (CABC)";                      B = 3; C = 0; D = 22.04;  E = 7.26;  F = "🟢 Achat";     G = "✅ Renforcer" },
    @{ Row = 23; A = "SAFCA CI (SAFC)";                         B = 2; C = 0; D = 14.58;  E = 7.24;  F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 24; A = "TRACTAFRIC MOTORS CI (PRSC)";             B = 2; C = 1; D = 11.99;  E = -2.78; F = "🟡 Observer"; G = "👀 À surveiller" },
    @{ Row = 25; A = "SETAO CI (STAC)";                         B = 2; C = 1; D = 7.39;   E = -7.2;  F = "🟡 Observer"; G = "👀 À surveiller" },
    @{ Row = 26; A = "ECOBANK TRANS. INCORP. TG (ETIT)";        B = 1; C = 0; D = 5.26;   E = 5.26;  F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 27; A = "BANK OF AFRICA SENEGAL (BOAS)";           B = 1; C = 0; D = 4.6;    E = 4.6;   F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 28; A = "LOTERIE NATIONALE DU BENIN (LNBB)";       B = 1; C = 1; D = 1.95;   E = 4.08;  F = "🟡 Observer"; G = "👀 À surveiller" },
    @{ Row = 29; A = "ONATEL BF (ONTBF)";                       B = 1; C = 2; D = 1.56;   E = -3.85; F = "🟡 Observer"; G = "👀 À surveiller" },
    @{ Row = 30; A = "SAPH CI (SPHC)";                          B = 1; C = 1; D = 0.64;   E = -6.81; F = "🟡 Observer"; G = "👀 À surveiller" },
    @{ Row = 31; A = "TOTAL";                                   B = 0; C = 3; D = 0;      E = 0;     F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 32; A = "ORAGROUP TOGO (ORGT)";                    B = 0; C = 1; D = -1.6;   E = -1.6;  F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 33; A = "NSIA BANQUE COTE D'IVOIRE (NSBC)";        B = 0; C = 1; D = -1.83;  E = -1.83; F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 34; A = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)";   B = 0; C = 1; D = -1.85;  E = -1.85; F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 35; A = "TOTALENERGIES MARKETING CI (TTLC)";       B = 0; C = 1; D = -2;     E = -2;    F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 36; A = "BERNABE CI (BNBC)";                       B = 0; C = 1; D = -3.16;  E = -3.16; F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 37; A = "NEI-CEDA CI (NEIC)";                      B = 0; C = 1; D = -3.65;  E = -3.65; F = "🟡 Observer"; G = "➖ Neutre" },
    @{ Row = 38; A = "FILTISAC CI (FTSC)";                      B = 1; C = 2; D = -5.84;  E = -5.75; F = "🟡 Observer"; G = "👀 À surveiller" }
)

foreach ($u in $titleUpdates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 7).Value = $u.G
}

# Rows 39-42 no longer appear in the refreshed ranking - remove them entirely
# (shifts the dimension down to A1:G38).
$ws.Range("A39:G42").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet "Top_YTD"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ytdUpdates = @(
    @{ Row = 2;  B = 523503.3 },
    @{ Row = 4;  B = 44078.75 },
    @{ Row = 5;  B = 35506.36 },
    @{ Row = 6;  B = 19793.9 },
    @{ Row = 7;  B = 10929.18 },
    @{ Row = 8;  B = 9830.120000000001 },
    @{ Row = 9;  B = 1734.45 },
    @{ Row = 10; B = 1244.25 },
    @{ Row = 11; B = 1219.61 }
)

foreach ($u in $ytdUpdates) {
    $ws2.Cells.Item($u.Row, 2).Value = $u.B
}
